# Cost.xlsx — "unify the conception of DataNode, DataTable, Entity"
#
# The sheet that described a generic "Property" table is renamed to
# "DataNode" to line up with the unified DataNode/DataTable/Entity naming
# used elsewhere in the project. The author's last on-sheet selection
# (before saving) also moved down to D39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Restore the author's last selection/active cell (D39) on the sheet
$ws.Range("D39").Select()
